$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Gift" -> "Savings Acc Interest" (amount + date updated)
$ws.Range("A2").Value = "Savings Acc Interest"
$ws.Range("B2").Value = 5026
$ws.Range("C2").Value = 45982.22928240741

# Row 3: "Salary Bonus" -> "FD Returns" (amount + date updated)
$ws.Range("A3").Value = "FD Returns"
$ws.Range("B3").Value = 15022
$ws.Range("C3").Value = 45981.22928240741

# Row 4: "Salary" stays, amount + date updated
$ws.Range("A4").Value = "Salary"
$ws.Range("B4").Value = 127000
$ws.Range("C4").Value = 45968.22928240741

# Row 5: brand new row - "Investment Returns"
$ws.Range("A5").Value = "Investment Returns"
$ws.Range("B5").Value = 10356
$ws.Range("C5").Value = 45966.22928240741

# Copy the date formatting (style index) from C4 onto the new C5 cell
# so it matches the existing date-formatted column exactly.
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
